$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.207.41'
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").Value = '1.788.62'
$ws.Range("E3").Value = '  -0.12%  '
$ws.Range("E4").Value = '  +0.12%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '226.09'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.39%  '
$ws.Range("E7").Value = '  +0.11%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '32.37'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("E10").Value = '  +0.08%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0948'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.88%  '
$ws.Range("D12").Value = '2.046.39'
$ws.Range("E12").Value = '  -0.10%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '11.11'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -3.01%  '
$ws.Range("D14").Value = '1.797.41'
$ws.Range("E14").Value = '  +0.34%  '
$ws.Range("E15").Value = '  +0.22%  '
$ws.Range("D16").Value = '34.185.74'
$ws.Range("E16").Value = '  +0.32%  '
$ws.Range("E17").Value = '  +0.26%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '67.98'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -0.17%  '
$ws.Range("D19").Value = '0.0₃0809'
$ws.Range("E19").Value = '  +3.32%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '246.17'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.88%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '11.01'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.58%  '
$ws.Range("E22").Value = '  +0.13%  '
$ws.Range("E23").Value = '  +1.60%  '
$ws.Range("E24").Value = '  +0.36%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '161.95'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("E26").Value = '  -0.47%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '16.34'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +0.26%  '
$ws.Range("E28").Value = '  +0.77%  '
$ws.Range("E29").Value = '  +0.36%  '
$ws.Range("E31").Value = '  -0.53%  '
$ws.Range("E32").Value = '  +2.33%  '
$ws.Range("E33").Value = '  +3.71%  '
$ws.Range("E34").Value = '  -2.24%  '
$ws.Range("D35").Value = '1.442.19'
$ws.Range("E35").Value = '  +1.93%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '2.56'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +8.12%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.665'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +2.75%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '1.05'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +1.17%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.0191'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.21%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '82.29'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +1.87%  '
$ws.Range("E41").Value = '  +1.38%  '
$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '13.89'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +2.57%  '
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.922'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -0.09%  '
$ws.Range("E44").Value = '  +0.96%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.0520'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +2.57%  '
$ws.Range("E46").Value = '  +0.18%  '
$ws.Range("E47").Value = '  +0.61%  '
$ws.Range("D48").Value = '1.946.08'
$ws.Range("E48").Value = '  -0.18%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '105.45'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -1.83%  '
$ws.Range("E50").Value = '  +0.14%  '
$ws.Range("E51").Value = '  -6.35%  '
